# Updated non-tracing input data
# The "N" row (B3/C3) used to hold hard-coded totals; replace them with
# formulas that sum the individual compartments in rows 16-22 so the
# total tracks the underlying (non-tracing) inputs automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ideal-format")

$ws.Range("B3").Formula = "=B16+B17+B18+B19+B20+B21+B22"
$ws.Range("C3").Formula = "=C16+C17+C18+C19+C20+C21+C22"

# Reflect where the author's cursor ended up after editing these cells.
$ws.Range("C3").Select()
